# 案件情報.xlsx — append scrape run 2025-10-07 01:16:07 (JST) on "ランサーズ" sheet.
#
# - refresh "取得日時" (A2:A17) to the new scrape timestamp
# - insert 2 brand-new listings at rows 11-12 (pushing the old rows 11-17 down to 13-19)
# - widen column B by one character (46 -> 47)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-07 01:16:07"

# --- widen column B (46 -> 47 chars); 46.17 is the ColumnWidth input that the
# --- engine's char-width -> OOXML-width rounding maps back onto exactly 47 ---
$ws.Columns.Item(2).ColumnWidth = 46.17

# --- refresh the scrape timestamp on all of the pre-existing rows (2-17) ---
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- make room for the 2 new listings: push old rows 11-17 down to 13-19 ---
$ws.Rows.Item(11).Resize(2).Insert()

# --- new row 11: iOS/Android app listing ---
$ws.Cells.Item(11, 1).Value = $newTimestamp
$ws.Cells.Item(11, 2).Value = "初回 iOSとAndroidのアプリ 課金(サブスク)"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5398382"
$ws.Cells.Item(11, 7).Value = 30
$ws.Cells.Item(11, 8).Value = "◇アプリ"

# --- new row 12: Ctrl+click folder app listing ---
$ws.Cells.Item(12, 1).Value = $newTimestamp
$ws.Cells.Item(12, 2).Value = "Ctrlを押しながらフォルダの上をクリックすると別窓で上の階層のフォルダが開くアプリの作成"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5408148"
$ws.Cells.Item(12, 7).Value = 30
$ws.Cells.Item(12, 8).Value = "◇アプリ"

# --- the 2 listings that fell off the bottom of the old range (old rows 16 &
# --- 17, now sitting at 18 & 19) need hyperlinks added on their URL cells ---
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), "https://www.lancers.jp/work/detail/5407390")
$ws.Cells.Item(18, 6).Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://www.lancers.jp/work/detail/5407516")
$ws.Cells.Item(19, 6).Style = "Hyperlink"
